# Apply the data-refresh edit described in the commit:
# "Update cfb_weather.xlsx with Timestamp 2025-09-14T16:23:15.489466"
#
# This updates:
#  - the shared Timestamp string used throughout sheet "FBS" (column AK)
#  - several "wind_dir_fg" (final-game wind direction) text cells on both
#    sheets, which were recomputed with the refreshed data
#  - a few numeric weather values (temp_fg / wind_fg / wind_diff) for the
#    first data row of sheet "Other"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FBS")
$ws2 = $wb.Worksheets.Item("Other")

# --- Timestamp column (AK) on sheet FBS: every data row shares the same
#     string value, so update each one to the new timestamp. ---
$newTimestamp = "2025-09-14T16:23:15.489466"
$lastRow = $ws1.Cells.Item($ws1.Rows.Count, 1).End(-4162).Row   # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $ws1.Range("AK$r").Value = $newTimestamp
}

# --- wind_dir_fg (column Q) updates on sheet FBS ---
$ws1.Range("Q5").Value  = "ENE"
$ws1.Range("Q7").Value  = "S"
$ws1.Range("Q10").Value = "NW"
$ws1.Range("Q11").Value = "NW"
$ws1.Range("Q13").Value = "S"
$ws1.Range("Q14").Value = "WSW"
$ws1.Range("Q17").Value = "SW"
$ws1.Range("Q22").Value = "WSW"
$ws1.Range("Q28").Value = "ENE"
$ws1.Range("Q29").Value = "NW"
$ws1.Range("Q33").Value = "WSW"
$ws1.Range("Q37").Value = "ENE"
$ws1.Range("Q39").Value = "S"
$ws1.Range("Q41").Value = "S"
$ws1.Range("Q43").Value = "WSW"
$ws1.Range("Q45").Value = "SSW"

# --- sheet Other, first data row: temp_fg / wind_fg / wind_diff ---
$ws2.Range("Q2").Value = 79.76000000000001
$ws2.Range("R2").Value = 9.1
$ws2.Range("W2").Value = -4.9

# --- wind_dir_fg (column S) updates on sheet Other ---
$ws2.Range("S14").Value = "SE"
$ws2.Range("S17").Value = "SW"
$ws2.Range("S23").Value = "S"
$ws2.Range("S25").Value = "ENE"
$ws2.Range("S27").Value = "SW"
$ws2.Range("S28").Value = "ENE"
$ws2.Range("S31").Value = "WSW"
$ws2.Range("S42").Value = "SSW"
$ws2.Range("S44").Value = "WSW"
$ws2.Range("S47").Value = "S"
$ws2.Range("S50").Value = "WSW"
